$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (Price / Volume(1h) columns).
# Values are written with a leading quote-prefix so Excel treats them as
# literal text (not numbers/dates), then the style is reset to "Normal"
# so no quotePrefix / text-number-format is left behind on the cell.
$updates = [ordered]@{
    'D2' = '26.430.00'
    'E2' = '  +1.02%  '
    'D3' = '1.723.92'
    'E3' = '  +1.87%  '
    'D4' = '1.001'
    'E4' = '  +0.30%  '
    'D5' = '244.49'
    'E5' = '  +2.17%  '
    'D7' = '0.4791'
    'E7' = '  +2.52%  '
    'D8' = '0.2678'
    'E8' = '  +1.69%  '
    'D9' = '0.06214'
    'E9' = '  +0.49%  '
    'D10' = '1.730.63'
    'E10' = '  +2.33%  '
    'D11' = '0.07113'
    'E11' = '  +0.89%  '
    'D12' = '15.68'
    'E12' = '  +3.42%  '
    'D13' = '0.6149'
    'E13' = '  +4.93%  '
    'D14' = '4.540'
    'E14' = '  +2.92%  '
    'E15' = '  +1.58%  '
    'D16' = '1.001'
    'E16' = '  +0.17%  '
    'D17' = '26.455.09'
    'E17' = '  +1.24%  '
    'E18' = '  +0.17%  '
    'D19' = '0.000006928'
    'E19' = '  +2.33%  '
    'D20' = '11.69'
    'E20' = '  +1.14%  '
    'D21' = '1.953.27'
    'D22' = '4.534'
    'E22' = '  +0.09%  '
    'D23' = '8.892'
    'E23' = '  +1.38%  '
    'D24' = '5.306'
    'E24' = '  +0.30%  '
    'D25' = '136.19'
    'E25' = '  +1.40%  '
    'D26' = '15.33'
    'E26' = '  +1.46%  '
    'E27' = '  +3.27%  '
    'D28' = '1.412'
    'E28' = '  +1.05%  '
    'D29' = '106.63'
    'E29' = '  +1.04%  '
    'D30' = '3.965'
    'E30' = '  -0.61%  '
    'D31' = '0.08032'
    'E31' = '  +3.56%  '
    'D32' = '3.730'
    'E32' = '  +1.56%  '
    'D33' = '0.04543'
    'E33' = '  +3.65%  '
    'D34' = '2.619'
    'E34' = '  +0.06%  '
    'D35' = '0.6356'
    'E35' = '  +3.28%  '
    'D36' = '0.9864'
    'E36' = '  +2.08%  '
    'D37' = '0.9350'
    'E37' = '  +1.21%  '
    'E38' = '  +5.59%  '
    'D39' = '107.43'
    'E39' = '  -2.92%  '
    'D40' = '2.403'
    'E40' = '  +0.86%  '
    'E41' = '  +0.60%  '
    'D42' = '0.01498'
    'E42' = '  +2.16%  '
    'D43' = '5.631'
    'E43' = '  +10.66%  '
    'D44' = '0.3903'
    'E44' = '  +3.60%  '
    'D45' = '6.967'
    'E45' = '  +12.12%  '
    'E46' = '  +5.20%  '
    'D47' = '0.05317'
    'E47' = '  -0.05%  '
    'E48' = '  +0.46%  '
    'D49' = '7.840'
    'E49' = '  +2.20%  '
    'D50' = '1.267'
    'E50' = '  +4.39%  '
    'D51' = '0.3416'
    'E51' = '  +2.03%  '
}

foreach ($ref in $updates.Keys) {
    $range = $ws.Range($ref)
    $range.Value = "'" + $updates[$ref]
    $range.Style = "Normal"
}
